$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 to make room for the new earliest forecast
# vector (shifts all existing data rows down by one, matching the new
# dimension A1:E53).
$ws.Rows.Item(2).Insert()

# The inserted row picks up a stray style (border/alignment copied from the
# row above); clear it so the new row matches the plain formatting used by
# every other data row, then restore the date display format on column A.
$ws.Range("A2:E2").ClearFormats()
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rewrite the full data block (A2:E53) with the corrected / regenerated
# forecast vectors (dates, y_0, y_0_forecast, y_1, y_1_forecast) in one shot.
$arr = New-Object 'object[,]' 52,5
$arr[0,0] = 39400
$arr[0,1] = 2007
$arr[0,2] = 3.145939949069287
$arr[0,3] = 2008
$arr[0,4] = 2.455972543253826
$arr[1,0] = 39583
$arr[1,1] = 2008
$arr[1,2] = 2.740959689118805
$arr[1,3] = 2009
$arr[1,4] = 3.206168778303486
$arr[2,0] = 39765
$arr[2,1] = 2008
$arr[2,2] = 1.769627576887389
$arr[2,3] = 2009
$arr[2,4] = 1.087227286828241
$arr[3,0] = 39948
$arr[3,1] = 2009
$arr[3,2] = -5.478010998490157
$arr[3,3] = 2010
$arr[3,4] = -2.013762956649334
$arr[4,0] = 40130
$arr[4,1] = 2009
$arr[4,2] = -4.774178217057779
$arr[4,3] = 2010
$arr[4,4] = -0.8523446516643385
$arr[5,0] = 40310
$arr[5,1] = 2010
$arr[5,2] = -0.1091898317121864
$arr[5,3] = 2011
$arr[5,4] = -1.305195642355683
$arr[6,0] = 40494
$arr[6,1] = 2010
$arr[6,2] = 1.97975191822708
$arr[6,3] = 2011
$arr[6,4] = 1.39052144387346
$arr[7,0] = 40676
$arr[7,1] = 2011
$arr[7,2] = 3.371423250978856
$arr[7,3] = 2012
$arr[7,4] = 0.8060632160631576
$arr[8,0] = 40862
$arr[8,1] = 2011
$arr[8,2] = 3.452886745653183
$arr[8,3] = 2012
$arr[8,4] = 2.798447799311043
$arr[9,0] = 41044
$arr[9,1] = 2012
$arr[9,2] = 1.627570629117536
$arr[9,3] = 2013
$arr[9,4] = 2.766358213445708
$arr[10,0] = 41228
$arr[10,1] = 2012
$arr[10,2] = 1.239479831392853
$arr[10,3] = 2013
$arr[10,4] = 2.047133666472267
$arr[11,0] = 41409
$arr[11,1] = 2013
$arr[11,2] = -0.03183655677961861
$arr[11,3] = 2014
$arr[11,4] = 1.102200073559878
$arr[12,0] = 41592
$arr[12,1] = 2013
$arr[12,2] = 0.2379616621361214
$arr[12,3] = 2014
$arr[12,4] = 1.062273708599726
$arr[13,0] = 41774
$arr[13,1] = 2014
$arr[13,2] = 1.812248956008733
$arr[13,3] = 2015
$arr[13,4] = 1.209672013646323
$arr[14,0] = 41957
$arr[14,1] = 2014
$arr[14,2] = 1.51977456621637
$arr[14,3] = 2015
$arr[14,4] = 0.6176326357196116
$arr[15,0] = 42137
$arr[15,1] = 2015
$arr[15,2] = 1.290465392296114
$arr[15,3] = 2016
$arr[15,4] = 0.9879295308886871
$arr[16,0] = 42321
$arr[16,1] = 2015
$arr[16,2] = 1.470039379455756
$arr[16,3] = 2016
$arr[16,4] = 1.577608035818301
$arr[17,0] = 42503
$arr[17,1] = 2016
$arr[17,2] = 1.57569012346459
$arr[17,3] = 2017
$arr[17,4] = 1.643656926428561
$arr[18,0] = 42689
$arr[18,1] = 2016
$arr[18,2] = 1.638797242243251
$arr[18,3] = 2017
$arr[18,4] = 1.369334405341616
$arr[19,0] = 42867
$arr[19,1] = 2017
$arr[19,2] = 1.73823635068906
$arr[19,3] = 2018
$arr[19,4] = 1.765380623247137
$arr[20,0] = 43053
$arr[20,1] = 2017
$arr[20,2] = 2.161565493242668
$arr[20,3] = 2018
$arr[20,4] = 2.486299099038347
$arr[21,0] = 43145
$arr[21,1] = 2018
$arr[21,2] = 2.50728418643813
$arr[21,3] = 2019
$arr[21,4] = 2.11231490846715
$arr[22,0] = 43235
$arr[22,1] = 2018
$arr[22,2] = 2.337818484846466
$arr[22,3] = 2019
$arr[22,4] = 2.076648015684435
$arr[23,0] = 43326
$arr[23,1] = 2018
$arr[23,2] = 2.354760705778203
$arr[23,3] = 2019
$arr[23,4] = 2.107524645430914
$arr[24,0] = 43418
$arr[24,1] = 2018
$arr[24,2] = 2.214251681313772
$arr[24,3] = 2019
$arr[24,4] = 1.362030665126834
$arr[25,0] = 43510
$arr[25,1] = 2019
$arr[25,2] = 0.8037559998091082
$arr[25,3] = 2020
$arr[25,4] = 1.649904670037805
$arr[26,0] = 43600
$arr[26,1] = 2019
$arr[26,2] = 0.8311911554373719
$arr[26,3] = 2020
$arr[26,4] = 1.758956425699298
$arr[27,0] = 43691
$arr[27,1] = 2019
$arr[27,2] = 0.7024402883234027
$arr[27,3] = 2020
$arr[27,4] = 1.302300993836147
$arr[28,0] = 43783
$arr[28,1] = 2019
$arr[28,2] = 0.6066442151010376
$arr[28,3] = 2020
$arr[28,4] = 0.7878236429522678
$arr[29,0] = 43875
$arr[29,1] = 2020
$arr[29,2] = 0.5367417164559685
$arr[29,3] = 2021
$arr[29,4] = 0.9721240557711397
$arr[30,0] = 43966
$arr[30,1] = 2020
$arr[30,2] = -1.538034740964334
$arr[30,3] = 2021
$arr[30,4] = -0.7351085756681308
$arr[31,0] = 44068
$arr[31,1] = 2020
$arr[31,2] = -4.65090747647452
$arr[31,3] = 2021
$arr[31,4] = -1.895157449498863
$arr[32,0] = 44159
$arr[32,1] = 2020
$arr[32,2] = -4.207901339433196
$arr[32,3] = 2021
$arr[32,4] = -0.4850133725290084
$arr[33,0] = 44251
$arr[33,1] = 2021
$arr[33,2] = 0.003696830084365388
$arr[33,3] = 2022
$arr[33,4] = -1.242205446257827
$arr[34,0] = 44341
$arr[34,1] = 2021
$arr[34,2] = 0.303920243687994
$arr[34,3] = 2022
$arr[34,4] = -0.6038293380915438
$arr[35,0] = 44432
$arr[35,1] = 2021
$arr[35,2] = 0.7583924418458787
$arr[35,3] = 2022
$arr[35,4] = 1.108416787477773
$arr[36,0] = 44525
$arr[36,1] = 2021
$arr[36,2] = 1.099928004397532
$arr[36,3] = 2022
$arr[36,4] = 1.102608990163567
$arr[37,0] = 44617
$arr[37,1] = 2022
$arr[37,2] = 1.880148611648913
$arr[37,3] = 2023
$arr[37,4] = 0.01066574587431646
$arr[38,0] = 44706
$arr[38,1] = 2022
$arr[38,2] = 1.91914784107321
$arr[38,3] = 2023
$arr[38,4] = 0.2384815980940092
$arr[39,0] = 44798
$arr[39,1] = 2022
$arr[39,2] = 2.236860175919531
$arr[39,3] = 2023
$arr[39,4] = 0.8831516962375607
$arr[40,0] = 44890
$arr[40,1] = 2022
$arr[40,2] = 2.310042359896225
$arr[40,3] = 2023
$arr[40,4] = 2.086256540666986
$arr[41,0] = 44981
$arr[41,1] = 2023
$arr[41,2] = 0.2425620590337463
$arr[41,3] = 2024
$arr[41,4] = 1.106935253696562
$arr[42,0] = 45071
$arr[42,1] = 2023
$arr[42,2] = -0.1211988132392205
$arr[42,3] = 2024
$arr[42,4] = 0.7707080878861294
$arr[43,0] = 45163
$arr[43,1] = 2023
$arr[43,2] = -0.09588622947416248
$arr[43,3] = 2024
$arr[43,4] = 0.8742015250004842
$arr[44,0] = 45254
$arr[44,1] = 2023
$arr[44,2] = 0.0464415346324687
$arr[44,3] = 2024
$arr[44,4] = 0.3491198177708599
$arr[45,0] = 45345
$arr[45,1] = 2024
$arr[45,2] = -0.297474409307219
$arr[45,3] = 2025
$arr[45,4] = 0.07171493608653101
$arr[46,0] = 45436
$arr[46,1] = 2024
$arr[46,2] = -0.04760886976447054
$arr[46,3] = 2025
$arr[46,4] = 0.5568966348730831
$arr[47,0] = 45534
$arr[47,1] = 2024
$arr[47,2] = -0.2385784141923808
$arr[47,3] = 2025
$arr[47,4] = -0.06923086958923186
$arr[48,0] = 45618
$arr[48,1] = 2024
$arr[48,2] = -0.3101476031197148
$arr[48,3] = 2025
$arr[48,4] = 0.2126457877301924
$arr[49,0] = 45713
$arr[49,1] = 2025
$arr[49,2] = -0.07765238411295838
$arr[49,3] = 2026
$arr[49,4] = -0.17084471036517
$arr[50,0] = 45800
$arr[50,1] = 2025
$arr[50,2] = 0.1245593350339691
$arr[50,3] = 2026
$arr[50,4] = 0.1722027100061974
$arr[51,0] = 45891
$arr[51,1] = 2025
$arr[51,2] = -0.0960403240804597
$arr[51,3] = 2026
$arr[51,4] = -0.1059622177528863

$ws.Range("A2:E53").Value = $arr
